# BDD Lunch & Learn - presentation mods
#
# 1. Add the "March 2, 2011" date line to the title slide's subtitle.
# 2. Swap the order of the last two slides ("For more info... / cukes.info"
#    and "Shawn Wallace / Andy Vida" contact info) so the contact-info
#    slide now comes second-to-last, and mark it hidden with a slow
#    2-second transition (it's kept around as a just-in-case backup slide).

$p = $ppt.ActivePresentation

# --- Title slide: add the date under the title ---
$titleSlide = $p.Slides.Item(1)
$titleSlide.Shapes.Item(2).TextFrame.TextRange.Text = "March 2, 2011"

# --- Reorder: move the "Shawn Wallace / Andy Vida" contact slide (11) so
#     it follows the "For more info.../cukes.info" slide (12) ---
$p.Slides.Item(11).MoveTo(12)

# --- Hide the (now relocated) contact-info slide and give it a slow
#     transition ---
$contactSlide = $p.Slides.Item(11)
$contactSlide.SlideShowTransition.Duration = 2
$contactSlide.SlideShowTransition.Speed = 1
$contactSlide.SlideShowTransition.Hidden = $true
